$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '29.312.50'
$ws.Range('E2').Value = '  +2.97%  '

# Row 3
$ws.Range('D3').Value = '1.895.48'
$ws.Range('E3').Value = '  +1.02%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.32%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.13%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.38%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5148'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.01%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3921'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.21%  '

# Row 9
$ws.Range('E9').Value = '  +0.37%  '

# Row 10
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.29'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.77%  '

# Row 11
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.115'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.67%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.269'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.51%  '

# Row 13
$ws.Range('D13').Value = '1.892.24'
$ws.Range('E13').Value = '  +0.83%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.69'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.07%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.276'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.43%  '

# Row 16
$ws.Range('E16').Value = '  -0.31%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '93.10'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.53%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001105'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.20%  '

# Row 19
$ws.Range('E19').Value = '  +0.37%  '

# Row 20
$ws.Range('E20').Value = '  +1.00%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.003'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.30%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.010'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.17%  '

# Row 23
$ws.Range('D23').Value = '29.326.79'

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.12'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.30%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.214'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.78%  '

# Row 26
$ws.Range('D26').Value = '2.108.06'
$ws.Range('E26').Value = '  +0.79%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '158.99'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.19%  '

# Row 28
$ws.Range('E28').Value = '  +1.21%  '

# Row 29
$ws.Range('E29').Value = '  +2.44%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.11'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.85%  '

# Row 31
$ws.Range('E31').Value = '  +1.30%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.144'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.51%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.660'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.59%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02480'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.04%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06576'
$ws.Range('D36').Style = 'Normal'

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2196'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.72%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '9.012'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.45%  '

# Row 39
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.199'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.65%  '

# Row 40
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.229'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.29%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6516'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.48%  '

# Row 42
$ws.Range('E42').Value = '  -2.07%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.26'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.12%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6057'
$ws.Range('D44').Style = 'Normal'

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.24'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.09%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.672'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.50%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.051'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.37%  '

# Row 48
$ws.Range('E48').Value = '  +2.13%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.56'
$ws.Range('D49').Style = 'Normal'

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.157'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.72%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '77.64'
$ws.Range('D51').Style = 'Normal'
